$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.342
$ws.Range("E4").Value = 16.424
$ws.Range("D7").Value = -7.066999999999998
$ws.Range("A8").Value = -22.218
$ws.Range("A10").Value = -21.856
$ws.Range("E11").Value = 17.109
$ws.Range("A12").Value = -21.589
$ws.Range("D14").Value = -7.802000000000001
$ws.Range("E14").Value = 17.257
$ws.Range("D15").Value = -8.196000000000002
$ws.Range("A18").Value = -22.166
$ws.Range("D18").Value = -8.4
$ws.Range("E18").Value = 16.366
$ws.Range("E19").Value = 16.452
$ws.Range("D20").Value = -7.367999999999999
$ws.Range("E21").Value = 16.421
$ws.Range("A25").Value = -21.691
$ws.Range("E27").Value = 16.38
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.048999999999999
$ws.Range("D31").Value = -7.613000000000001
$ws.Range("E31").Value = 16.822
$ws.Range("D35").Value = -7.737
$ws.Range("A37").Value = -20.105
$ws.Range("E38").Value = 16.591
$ws.Range("D40").Value = -7.587999999999999
$ws.Range("E42").Value = 16.406
$ws.Range("D44").Value = -7.421000000000001
$ws.Range("E44").Value = 16.937
$ws.Range("E47").Value = 16.466
$ws.Range("D50").Value = -8.104999999999999
$ws.Range("D54").Value = -8.100000000000001
$ws.Range("A55").Value = -22.311
$ws.Range("E56").Value = 16.289
$ws.Range("E58").Value = 16.448
$ws.Range("E65").Value = 17.301
$ws.Range("A68").Value = -21.534
$ws.Range("D68").Value = -6.778
$ws.Range("E73").Value = 16.402
$ws.Range("D76").Value = -7.672
$ws.Range("A77").Value = -20.938
$ws.Range("A78").Value = -20.22
$ws.Range("A79").Value = -21.798
$ws.Range("A80").Value = -20.203
$ws.Range("A81").Value = -21.782
$ws.Range("A82").Value = -22.261
$ws.Range("A84").Value = -22.106
$ws.Range("D87").Value = -8.297000000000001
$ws.Range("D88").Value = -8.259000000000002
$ws.Range("E90").Value = 16.415
$ws.Range("D92").Value = -7.306
$ws.Range("E92").Value = 17.069
$ws.Range("E94").Value = 17.686
$ws.Range("E95").Value = 17.121
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.404
$ws.Range("A101").Value = -21.41
$ws.Range("D101").Value = -7.712999999999999
$ws.Range("E101").Value = 16.733
$ws.Range("A102").Value = -20.163
$ws.Range("D102").Value = -8.090999999999999
